# Weekly update: a new price record for "Arándano (blue)" at Macroferia
# Regional de Talca is inserted as the new first row of the data block
# (row 57), pushing every subsequent record down by one row. The sheet's
# last existing record (previously row 148) lands on the newly created
# row 149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 57 - this shifts rows 57:148 down to 58:149
# and extends the sheet's used range/dimension accordingly.
$ws.Rows.Item(57).Insert()

# Populate the new row with this week's record.
$ws.Cells.Item(57, 1).Value  = 5
$ws.Cells.Item(57, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(57, 3).Value  = "Maule"
$ws.Cells.Item(57, 4).Value  = 45259
$ws.Cells.Item(57, 5).Value  = 7
$ws.Cells.Item(57, 6).Value  = "Fruta"
$ws.Cells.Item(57, 7).Value  = 100101
$ws.Cells.Item(57, 8).Value  = "Berries"
$ws.Cells.Item(57, 9).Value  = 100101001
$ws.Cells.Item(57, 10).Value = "Arándano (blue)"
$ws.Cells.Item(57, 11).Value = "Sin especificar"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 50
$ws.Cells.Item(57, 14).Value = 5000
$ws.Cells.Item(57, 15).Value = 5000
$ws.Cells.Item(57, 16).Value = 5000
$ws.Cells.Item(57, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(57, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(57, 19).Value = 2500
$ws.Cells.Item(57, 20).Value = 2
